# Regenerate save_data column G ("K") using the new strikeout calculation
# (replacing the old Strike# derived values) and write the recalculated values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2 = 1
    3 = 0
    4 = 1
    5 = 1
    6 = 0
    7 = 5
    8 = 2
    9 = 5
    10 = 1
    11 = 4
    12 = 4
    13 = 2
    14 = 7
    15 = 1
    16 = 1
    17 = 6
    18 = 2
    19 = 1
    20 = 4
    21 = 4
    22 = 3
    23 = 1
    24 = 2
    25 = 2
    26 = 5
    27 = 4
    28 = 1
    30 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
